# "game mechanics done, information"
# Update the thesis page-count tracker: Theory chapter written pages
# increased 8 -> 12, and the Discussion day-tally (I13) increased 0 -> 4.
# Dependent totals (B11, D11, B12, I31) and the day-counter ratios
# (F2, F3, F4 - driven by the volatile TODAY() formula) recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Theory: Written pages 8 -> 12
$ws.Range("B4").Value = 12

# Discussion row tally 0 -> 4
$ws.Range("I13").Value = 4

# Move the active selection to C20, matching the author's last cursor
# position when they saved the workbook.
$ws.Activate()
$ws.Range("C20").Select()
